$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.8775636666666666
$ws.Range("H2").Value = 2.632691
$ws.Range("I2").Value = 0.1887436506618166
$ws.Range("J2").Value = 0.2083714858314108
$ws.Range("M2").Value = 5.978421000000001
$ws.Range("N2").Value = 17.935263
$ws.Range("O2").Value = 0.05704457007880161
$ws.Range("P2").Value = 0.06242884486533885
$ws.Range("Q2").Value = 5.246445053637
$ws.Range("R2").Value = 47.21800548273301
$ws.Range("S2").Value = 0.01076680040710684
$ws.Range("T2").Value = 0.01300839116332929
# Row 3
$ws.Range("G3").Value = 0.8775636666666666
$ws.Range("H3").Value = 2.632691
$ws.Range("I3").Value = 0.1887436506618166
$ws.Range("J3").Value = 0.2083714858314108
$ws.Range("O3").Value = 0.6646576013185088
$ws.Range("P3").Value = 0.7273927426214574
$ws.Range("Q3").Value = 61.12921142157167
$ws.Range("R3").Value = 550.162902794145
$ws.Range("S3").Value = 0.1254499021129816
$ws.Range("T3").Value = 0.151567906563018
# Row 4
$ws.Range("G4").Value = 0.8775636666666666
$ws.Range("H4").Value = 2.632691
$ws.Range("I4").Value = 0.1887436506618166
$ws.Range("J4").Value = 0.2083714858314108
$ws.Range("M4").Value = 1.290243
$ws.Range("N4").Value = 3.870729
$ws.Range("O4").Value = 0.01231116999491725
$ws.Range("P4").Value = 0.01347318632889677
$ws.Range("Q4").Value = 1.132270377971
$ws.Range("R4").Value = 10.190433401739
$ws.Range("S4").Value = 0.002323655168758901
$ws.Range("T4").Value = 0.002807427854235671
# Row 5
$ws.Range("G5").Value = 0.8775636666666666
$ws.Range("H5").Value = 2.632691
$ws.Range("I5").Value = 0.1887436506618166
$ws.Range("J5").Value = 0.2083714858314108
$ws.Range("M5").Value = 27.1166075
$ws.Range("N5").Value = 54.233215
$ws.Range("O5").Value = 0.2587397603536297
$ws.Range("P5").Value = 0.1887743138075849
$ws.Range("Q5").Value = 23.79654950526083
$ws.Range("R5").Value = 142.779297031565
$ws.Range("S5").Value = 0.04883548694050762
$ws.Range("T5").Value = 0.03933518425489148
# Row 6
$ws.Range("G6").Value = 0.8775636666666666
$ws.Range("H6").Value = 2.632691
$ws.Range("I6").Value = 0.1887436506618166
$ws.Range("J6").Value = 0.2083714858314108
$ws.Range("M6").Value = 0.759494
$ws.Range("N6").Value = 2.278482
$ws.Range("O6").Value = 0.00724689825414258
$ws.Range("P6").Value = 0.007930912376722157
$ws.Range("Q6").Value = 0.6665043394513334
$ws.Range("R6").Value = 5.998539055061999
$ws.Range("S6").Value = 0.001367806032461616
$ws.Range("T6").Value = 0.001652575995936321
# Row 7
$ws.Range("G7").Value = 2.458038666666667
$ws.Range("H7").Value = 7.374116000000001
$ws.Range("I7").Value = 0.5286672739959656
$ws.Range("J7").Value = 0.5836444564186148
$ws.Range("M7").Value = 5.978421000000001
$ws.Range("N7").Value = 17.935263
$ws.Range("O7").Value = 0.05704457007880161
$ws.Range("P7").Value = 0.06242884486533885
$ws.Range("Q7").Value = 14.695189983612
$ws.Range("R7").Value = 132.256709852508
$ws.Range("S7").Value = 0.03015759735983187
$ws.Range("T7").Value = 0.03643624922627273
# Row 8
$ws.Range("G8").Value = 2.458038666666667
$ws.Range("H8").Value = 7.374116000000001
$ws.Range("I8").Value = 0.5286672739959656
$ws.Range("J8").Value = 0.5836444564186148
$ws.Range("O8").Value = 0.6646576013185088
$ws.Range("P8").Value = 0.7273927426214574
$ws.Range("Q8").Value = 171.2217256074467
$ws.Range("S8").Value = 0.3513827222297534
$ws.Range("T8").Value = 0.4245387418701459
# Row 9
$ws.Range("G9").Value = 2.458038666666667
$ws.Range("H9").Value = 7.374116000000001
$ws.Range("I9").Value = 0.5286672739959656
$ws.Range("J9").Value = 0.5836444564186148
$ws.Range("M9").Value = 1.290243
$ws.Range("N9").Value = 3.870729
$ws.Range("O9").Value = 0.01231116999491725
$ws.Range("P9").Value = 0.01347318632889677
$ws.Range("Q9").Value = 3.171467183396
$ws.Range("R9").Value = 28.543204650564
$ws.Range("S9").Value = 0.006508512680913831
$ws.Range("T9").Value = 0.00786355051115567
# Row 10
$ws.Range("G10").Value = 2.458038666666667
$ws.Range("H10").Value = 7.374116000000001
$ws.Range("I10").Value = 0.5286672739959656
$ws.Range("J10").Value = 0.5836444564186148
$ws.Range("M10").Value = 27.1166075
$ws.Range("N10").Value = 54.233215
$ws.Range("O10").Value = 0.2587397603536297
$ws.Range("P10").Value = 0.1887743138075849
$ws.Range("Q10").Value = 66.65366974382334
$ws.Range("R10").Value = 399.9220184629401
$ws.Range("S10").Value = 0.1367872437805228
$ws.Range("T10").Value = 0.1101770817680249
# Row 11
$ws.Range("G11").Value = 2.458038666666667
$ws.Range("H11").Value = 7.374116000000001
$ws.Range("I11").Value = 0.5286672739959656
$ws.Range("J11").Value = 0.5836444564186148
$ws.Range("M11").Value = 0.759494
$ws.Range("N11").Value = 2.278482
$ws.Range("O11").Value = 0.00724689825414258
$ws.Range("P11").Value = 0.007930912376722157
$ws.Range("Q11").Value = 1.866865619101334
$ws.Range("R11").Value = 16.801790571912
$ws.Range("S11").Value = 0.00383119794494368
$ws.Range("T11").Value = 0.004628833043015668
# Row 12
$ws.Range("G12").Value = 1.313898
$ws.Range("H12").Value = 2.627796
$ws.Range("I12").Value = 0.2825890753422177
$ws.Range("J12").Value = 0.2079840577499744
$ws.Range("M12").Value = 5.978421000000001
$ws.Range("N12").Value = 17.935263
$ws.Range("O12").Value = 0.05704457007880161
$ws.Range("P12").Value = 0.06242884486533885
$ws.Range("Q12").Value = 7.855035395058001
$ws.Range("R12").Value = 47.130212370348
$ws.Range("S12").Value = 0.01612017231186288
$ws.Range("T12").Value = 0.01298420447573683
# Row 13
$ws.Range("G13").Value = 1.313898
$ws.Range("H13").Value = 2.627796
$ws.Range("I13").Value = 0.2825890753422177
$ws.Range("J13").Value = 0.2079840577499744
$ws.Range("O13").Value = 0.6646576013185088
$ws.Range("P13").Value = 0.7273927426214574
$ws.Range("Q13").Value = 91.52332950777
$ws.Range("R13").Value = 549.1399770466199
$ws.Range("S13").Value = 0.1878249769757738
$ws.Range("T13").Value = 0.1512860941882935
# Row 14
$ws.Range("G14").Value = 1.313898
$ws.Range("H14").Value = 2.627796
$ws.Range("I14").Value = 0.2825890753422177
$ws.Range("J14").Value = 0.2079840577499744
$ws.Range("M14").Value = 1.290243
$ws.Range("N14").Value = 3.870729
$ws.Range("O14").Value = 0.01231116999491725
$ws.Range("P14").Value = 0.01347318632889677
$ws.Range("Q14").Value = 1.695247697214
$ws.Range("R14").Value = 10.171486183284
$ws.Range("S14").Value = 0.003479002145244522
$ws.Range("T14").Value = 0.002802207963505433
# Row 15
$ws.Range("G15").Value = 1.313898
$ws.Range("H15").Value = 2.627796
$ws.Range("I15").Value = 0.2825890753422177
$ws.Range("J15").Value = 0.2079840577499744
$ws.Range("M15").Value = 27.1166075
$ws.Range("N15").Value = 54.233215
$ws.Range("O15").Value = 0.2587397603536297
$ws.Range("P15").Value = 0.1887743138075849
$ws.Range("Q15").Value = 35.628456361035
$ws.Range("R15").Value = 142.51382544414
$ws.Range("S15").Value = 0.07311702963259921
$ws.Range("T15").Value = 0.03926204778466854
# Row 16
$ws.Range("G16").Value = 1.313898
$ws.Range("H16").Value = 2.627796
$ws.Range("I16").Value = 0.2825890753422177
$ws.Range("J16").Value = 0.2079840577499744
$ws.Range("M16").Value = 0.759494
$ws.Range("N16").Value = 2.278482
$ws.Range("O16").Value = 0.00724689825414258
$ws.Range("P16").Value = 0.007930912376722157
$ws.Range("Q16").Value = 0.9978976476120001
$ws.Range("R16").Value = 5.987385885671999
$ws.Range("S16").Value = 0.002047894276737283
$ws.Range("T16").Value = 0.001649503337770168
